$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Q3"
$ws.Range("D1").Value = "Q4"

$ws.Range("A2:A19").Interior.Color = 65535
$ws.Range("B2:B19").Interior.ThemeColor = 3
$ws.Range("C2:C19").Interior.Color = 5296274
$ws.Range("D2:D19").Interior.Color = 6299648
